# Update the accelerometer sample window: drop the first 11 rows of the
# previous window, keep the remaining 9, and append 21 new rows so the
# sheet now spans A1:C31 (30 data rows under the x/y/z header).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(8.708540916442871, 0.2552618980407715, 13.9090576171875),
    @(24.69548416137696, -32.37852478027344, -13.10801219940186),
    @(9.214871406555176, -8.127243041992188, -7.553183078765869),
    @(-1.296450614929199, -25.33206176757812, 4.173199653625488),
    @(3.13936185836792, 1.159258842468261, 10.52872467041016),
    @(-5.014933586120605, 3.021981239318848, 2.716179370880127),
    @(12.84034442901611, -46.83102416992188, 6.994163036346436),
    @(39.45174789428711, 18.51260566711425, -32.22710418701172),
    @(-25.77109146118164, -4.628503799438477, -10.99930953979492),
    @(-25.15040969848633, -52.22837066650391, 40.054931640625),
    @(-0.4132165908813476, 5.698812484741211, 7.912927150726318),
    @(17.88103866577148, 1.754395723342896, 14.25908279418945),
    @(-0.3170166015625, -36.44598770141602, 3.271512508392334),
    @(10.27998065948486, 5.400233268737793, -10.95361423492432),
    @(-23.53899955749512, -6.704813957214356, -13.04527854919434),
    @(17.90056037902832, -19.72356796264648, -5.512706756591797),
    @(-16.73113822937012, 13.58313274383545, 5.493848323822021),
    @(4.27336597442627, -17.7934684753418, 25.06283569335937),
    @(-30.59529876708984, -20.34370040893555, -32.5362548828125),
    @(4.83027458190918, 5.810257911682129, -21.42179489135743),
    @(-11.62696266174316, -21.11956596374512, -16.67885398864746),
    @(-13.2412109375, 3.214882850646973, 22.37087440490723),
    @(-15.93152809143066, 19.48006248474121, 7.360580921173096),
    @(32.82965087890625, -70.24890899658203, 12.97732639312744),
    @(-13.05622291564941, 5.201236724853516, -37.13468551635742),
    @(-35.03109741210938, -2.42187762260437, -11.56947135925293),
    @(-34.99357986450195, -38.73427581787109, 22.45726585388184),
    @(-5.411000728607178, 11.04445457458496, 9.49652862548828),
    @(0.176581859588623, 12.73018646240234, 1.936479568481445),
    @(9.547159194946287, -75.66817474365234, 32.35793304443359),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}